$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact literal formatting (avoid numeric auto-conversion
# of values like "0.630" or "0.0804" into floating point numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.883.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.979.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.98"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.91%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0804"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.77%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.78"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.848"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.268.83"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.985.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.781.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.26%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.76"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.51"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +19.26%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.44%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.81%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -9.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.24"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.368.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.29"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.97"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.79%  "
